$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1930
$ws.Range("I29").Value = 1950
$ws.Range("J29").Value = 1890
$ws.Range("K29").Value = 5850
$ws.Range("L29").Value = 5670
$ws.Range("M29").Value = -5569
$ws.Range("N29").Value = -6232
$ws.Range("H38").Value = 1969.75
$ws.Range("J38").Value = 3745
$ws.Range("L38").Value = 11235
$ws.Range("N38").Value = -11979
$ws.Range("H58").Value = 5277.9165
$ws.Range("I58").Value = 803.1667
$ws.Range("J58").Value = 9752.666999999999
$ws.Range("K58").Value = 2409.5001
$ws.Range("L58").Value = 29258.001
$ws.Range("M58").Value = -2259.5001
$ws.Range("N58").Value = -29558.001
$ws.Range("H62").Value = 3999.6667
$ws.Range("I62").Value = 4499.5
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 4499.5
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -3875.5
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 3999.6667
$ws.Range("I65").Value = 4499.5
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 22497.5
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -19377.5
$ws.Range("N65").Value = -21240
$ws.Range("H112").Value = 1108.3334
$ws.Range("J112").Value = 1108.3334
$ws.Range("L112").Value = 3325.0002
$ws.Range("N112").Value = -5541.0002
$ws.Range("H132").Value = 1371.1428
$ws.Range("I132").Value = 1266.5
$ws.Range("K132").Value = 3799.5
$ws.Range("M132").Value = -1269.5
$ws.Range("H135").Value = 3609.5715
$ws.Range("I135").Value = 4566.2
$ws.Range("K135").Value = 41095.8
$ws.Range("M135").Value = -38560.8
$ws.Range("H138").Value = 3896.5122
$ws.Range("I138").Value = 2715.0833
$ws.Range("J138").Value = 4385.3794
$ws.Range("K138").Value = 8145.249899999999
$ws.Range("L138").Value = 13156.1382
$ws.Range("M138").Value = -3005.249899999999
$ws.Range("N138").Value = -23436.1382

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2417.5881
$ws.Range("I45").Value = 2165.5
$ws.Range("J45").Value = 2777.7144
$ws.Range("K45").Value = 2165.5
$ws.Range("L45").Value = 2777.7144
$ws.Range("M45").Value = -1788.5
$ws.Range("N45").Value = -3531.7144
$ws.Range("H61").Value = 1887.7858
$ws.Range("I61").Value = 1754.52
$ws.Range("K61").Value = 1754.52
$ws.Range("M61").Value = -1542.52
$ws.Range("H122").Value = 10800
$ws.Range("I122").Value = 10800
$ws.Range("K122").Value = 32400
$ws.Range("M122").Value = -29950
$ws.Range("H136").Value = 1887.7858
$ws.Range("I136").Value = 1754.52
$ws.Range("K136").Value = 5263.559999999999
$ws.Range("M136").Value = -2713.559999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H107").Value = 1930.5
$ws.Range("I107").Value = 2240.6667
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 2240.6667
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -320.6667000000002
$ws.Range("N107").Value = -4840

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4349.12
$ws.Range("I31").Value = 2545.0667
$ws.Range("K31").Value = 2545.0667
$ws.Range("M31").Value = -2250.0667
$ws.Range("H34").Value = 4349.12
$ws.Range("I34").Value = 2545.0667
$ws.Range("K34").Value = 2545.0667
$ws.Range("M34").Value = -2343.0667
$ws.Range("H58").Value = 2285.1667
$ws.Range("I58").Value = 2285.1667
$ws.Range("K58").Value = 2285.1667
$ws.Range("M58").Value = -2082.1667
$ws.Range("H99").Value = 1497.3334
$ws.Range("I99").Value = 1683.7142
$ws.Range("K99").Value = 1683.7142
$ws.Range("M99").Value = -185.7141999999999
$ws.Range("H107").Value = 853.36365
$ws.Range("I107").Value = 504.77777
$ws.Range("K107").Value = 504.77777
$ws.Range("M107").Value = 1415.22223
$ws.Range("H126").Value = 1497.3334
$ws.Range("I126").Value = 1683.7142
$ws.Range("K126").Value = 5051.142599999999
$ws.Range("M126").Value = -2581.142599999999
$ws.Range("H135").Value = 124549.5
$ws.Range("J135").Value = 124549.5
$ws.Range("L135").Value = 124549.5
$ws.Range("N135").Value = -134689.5
$ws.Range("H136").Value = 2285.1667
$ws.Range("I136").Value = 2285.1667
$ws.Range("K136").Value = 6855.500100000001
$ws.Range("M136").Value = -4305.500100000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 860
$ws.Range("I5").Value = 731.5
$ws.Range("J5").Value = 988.5
$ws.Range("K5").Value = 2194.5
$ws.Range("L5").Value = 2965.5
$ws.Range("M5").Value = -2082.5
$ws.Range("N5").Value = -3189.5
$ws.Range("H98").Value = 4042.4285
$ws.Range("J98").Value = 4566.1665
$ws.Range("L98").Value = 13698.4995
$ws.Range("N98").Value = -16694.4995
$ws.Range("H114").Value = 4468.2
$ws.Range("I114").Value = 3085.5
$ws.Range("J114").Value = 9999
$ws.Range("K114").Value = 9256.5
$ws.Range("L114").Value = 29997
$ws.Range("M114").Value = -6002.5
$ws.Range("N114").Value = -36505
$ws.Range("H132").Value = 4495.6665
$ws.Range("I132").Value = 2665.3333
$ws.Range("K132").Value = 23987.9997
$ws.Range("M132").Value = -21457.9997
$ws.Range("H135").Value = 860
$ws.Range("I135").Value = 731.5
$ws.Range("J135").Value = 988.5
$ws.Range("K135").Value = 6583.5
$ws.Range("L135").Value = 8896.5
$ws.Range("M135").Value = -4048.5
$ws.Range("N135").Value = -13966.5
$ws.Range("H138").Value = 2680
$ws.Range("I138").Value = 2680
$ws.Range("K138").Value = 8040
$ws.Range("M138").Value = -2900

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 35081
$ws.Range("I62").Value = 10077
$ws.Range("K62").Value = 10077
$ws.Range("M62").Value = -9391
$ws.Range("H65").Value = 35081
$ws.Range("I65").Value = 10077
$ws.Range("K65").Value = 30231
$ws.Range("M65").Value = -26799
$ws.Range("H97").Value = 853.125
$ws.Range("I97").Value = 689.2857
$ws.Range("K97").Value = 689.2857
$ws.Range("M97").Value = -193.2857
$ws.Range("H130").Value = 63332.5
$ws.Range("J130").Value = 63332.5
$ws.Range("L130").Value = 63332.5
$ws.Range("N130").Value = -73372.5
$ws.Range("H132").Value = 3212.524
$ws.Range("I132").Value = 2748.1667
$ws.Range("K132").Value = 8244.500100000001
$ws.Range("M132").Value = -5714.500100000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H68").Value = 2999
$ws.Range("I68").Value = 2999
$ws.Range("K68").Value = 2999
$ws.Range("M68").Value = -2250
$ws.Range("H71").Value = 2999
$ws.Range("I71").Value = 2999
$ws.Range("K71").Value = 14995
$ws.Range("M71").Value = -11251
$ws.Range("H82").Value = 300
$ws.Range("I82").Value = 300
$ws.Range("K82").Value = 300
$ws.Range("M82").Value = 61
$ws.Range("H85").Value = 300
$ws.Range("I85").Value = 300
$ws.Range("K85").Value = 300
$ws.Range("M85").Value = 948
$ws.Range("H132").Value = 4386.3125
$ws.Range("I132").Value = 4387.25
$ws.Range("J132").Value = 4385.375
$ws.Range("K132").Value = 13161.75
$ws.Range("L132").Value = 13156.125
$ws.Range("M132").Value = -10631.75
$ws.Range("N132").Value = -18216.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2541.8333
$ws.Range("I81").Value = 2541.8333
$ws.Range("K81").Value = 5083.6666
$ws.Range("M81").Value = -4022.6666
$ws.Range("H84").Value = 2541.8333
$ws.Range("I84").Value = 2541.8333
$ws.Range("K84").Value = 25418.333
$ws.Range("M84").Value = -20114.333
$ws.Range("H122").Value = 2166.3333
$ws.Range("I122").Value = 2166.3333
$ws.Range("K122").Value = 6498.999899999999
$ws.Range("M122").Value = -4048.999899999999
